$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- The re-save normalized the table's font (the bold header font and the
# 10pt body font were dropped back to the workbook's plain default font) ---
$table = $ws.Range("A1:D9")
$table.Font.Bold = $false
$table.Font.Size = 11

# --- Core content edit: fill in the "200m - F" row (row 5) with the 2024
# Olympics results (Winner / 2nd place / 3rd place) ---
# The freshly-typed "Gabrielle Thomas" cell picked up a plain black font with
# no border (rather than inheriting the row's bordered/default style), so it
# is entered first and reformatted before the other two cells are filled in.
$ws.Range("B5").ClearFormats()
$ws.Range("B5").Value = "Gabrielle Thomas"
$ws.Range("B5").Font.Color = 0

$ws.Range("C5").Value = "Julien Alfred"
$ws.Range("D5").Value = "Brittany Brown"

# --- Move the active selection, matching the cursor position at save time ---
[void]$ws.Range("D20").Select()
